$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.174.89'
$ws.Range("E2").Value = '  -0.98%  '

$ws.Range("D3").Value = '2.273.95'
$ws.Range("E3").Value = '  -1.37%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''299.65'
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("D6").Value = '''95.59'
$ws.Range("E6").Value = '  -3.82%  '

$ws.Range("D7").Value = '''0.495'
$ws.Range("E7").Value = '  -2.24%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -2.44%  '

$ws.Range("D10").Value = '''33.08'
$ws.Range("E10").Value = '  -4.70%  '

$ws.Range("D11").Value = '''0.0788'
$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("D12").Value = '''48.34'
$ws.Range("E12").Value = '  -6.54%  '

$ws.Range("E13").Value = '  +1.37%  '

$ws.Range("D14").Value = '''16.05'
$ws.Range("E14").Value = '  +2.42%  '

$ws.Range("D15").Value = '''6.68'
$ws.Range("E15").Value = '  -0.75%  '

$ws.Range("D16").Value = '2.626.41'
$ws.Range("E16").Value = '  -1.27%  '

$ws.Range("D17").Value = '2.279.96'

$ws.Range("D18").Value = '''0.787'
$ws.Range("E18").Value = '  -2.11%  '

$ws.Range("D19").Value = '42.123.73'
$ws.Range("E19").Value = '  -0.93%  '

$ws.Range("D20").Value = '''11.69'
$ws.Range("E20").Value = '  +2.02%  '

$ws.Range("D21").Value = '0.0₃0891'
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = '''5.99'
$ws.Range("E22").Value = '  -1.39%  '

$ws.Range("D23").Value = '''66.30'
$ws.Range("E23").Value = '  -2.22%  '

$ws.Range("D24").Value = '''235.46'
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("D27").Value = '''2.46'
$ws.Range("E27").Value = '  -2.21%  '

$ws.Range("D28").Value = '''23.81'
$ws.Range("E28").Value = '  -4.50%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '''167.61'
$ws.Range("E29").Value = '  +2.03%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''2.07'
$ws.Range("E30").Value = '  -4.85%  '

$ws.Range("D31").Value = '''33.69'
$ws.Range("E31").Value = '  -2.78%  '

$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("E34").Value = '  +5.98%  '

$ws.Range("E35").Value = '  -2.24%  '

$ws.Range("D36").Value = '''16.75'
$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("D37").Value = '''2.34'
$ws.Range("E37").Value = '  -3.07%  '

$ws.Range("D38").Value = '''0.0689'
$ws.Range("E38").Value = '  -1.94%  '

$ws.Range("E39").Value = '  -2.86%  '

$ws.Range("E40").Value = '  -1.62%  '

$ws.Range("E41").Value = '  -2.21%  '

$ws.Range("E42").Value = '  -4.33%  '

$ws.Range("D43").Value = '''2.31'
$ws.Range("E43").Value = '  -5.89%  '

$ws.Range("D44").Value = '1.959.76'
$ws.Range("E44").Value = '  -0.42%  '

$ws.Range("D45").Value = '''0.0277'
$ws.Range("E45").Value = '  -1.20%  '

$ws.Range("D46").Value = '''17.60'
$ws.Range("E46").Value = '  -4.39%  '

$ws.Range("D47").Value = '''9.57'
$ws.Range("E47").Value = '  -6.27%  '

$ws.Range("D48").Value = '''2.78'
$ws.Range("E48").Value = '  -3.80%  '

$ws.Range("D49").Value = '2.496.69'
$ws.Range("E49").Value = '  -1.29%  '

$ws.Range("D50").Value = '''52.34'
$ws.Range("E50").Value = '  -5.86%  '

$ws.Range("E51").Value = '  -4.45%  '
